# Update column F ("dSF") values on the active sheet to reflect the
# re-pulled data / mean-calculation fix described in the commit message.
# Only column F values change; every other column is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 1
    4  = -2
    5  = 1
    6  = 4
    7  = -4
    8  = 3
    9  = -2
    11 = 1
    12 = 10
    13 = -1
    14 = -2
    15 = 3
    16 = 2
    17 = -2
    18 = 6
    19 = 5
    20 = -2
    21 = 2
    22 = 4
    23 = -2
    24 = -2
    25 = -3
    26 = 3
    27 = -1
    28 = 3
    29 = -2
    30 = 3
    31 = 5
    32 = -3
    33 = -1
    34 = 9
    36 = -1
    38 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
